$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 4.4
$ws.Range("Q2").Value = 1.78
$ws.Range("T2").Value = 1.68
$ws.Range("U2").Value = 2.36
$ws.Range("AN2").Value = 29
$ws.Range("AO2").Value = 15

# Row 3
$ws.Range("G3").Value = 2.06
$ws.Range("J3").Value = 3.6
$ws.Range("N3").Value = 1.02
$ws.Range("S3").Value = 1.01
$ws.Range("W3").Value = 1.94

# Row 4
$ws.Range("S4").Value = 2.06
$ws.Range("Y4").Value = 980
$ws.Range("Z4").Value = 80

# Row 5
$ws.Range("F5").Value = 2.14
$ws.Range("G5").Value = 2.48
$ws.Range("H5").Value = 3.35
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 3
$ws.Range("N5").Value = 1.01
$ws.Range("O5").Value = 1.42
$ws.Range("P5").Value = 1.66
$ws.Range("Q5").Value = 2.14
$ws.Range("R5").Value = 1.19
$ws.Range("S5").Value = 3.2
$ws.Range("W5").Value = 1.67
$ws.Range("X5").Value = 16
$ws.Range("Y5").Value = 17.5
$ws.Range("Z5").Value = 38
$ws.Range("AA5").Value = 100
$ws.Range("AB5").Value = 12.5
$ws.Range("AC5").Value = 10.5
$ws.Range("AD5").Value = 23
$ws.Range("AE5").Value = 75
$ws.Range("AF5").Value = 19
$ws.Range("AG5").Value = 16
$ws.Range("AH5").Value = 29
$ws.Range("AI5").Value = 95
$ws.Range("AJ5").Value = 44
$ws.Range("AK5").Value = 38
$ws.Range("AL5").Value = 65

# Row 6
$ws.Range("G6").Value = 1.95
$ws.Range("H6").Value = 4.7
$ws.Range("N6").Value = 1.62
$ws.Range("O6").Value = 1.42
$ws.Range("R6").Value = 1.18
$ws.Range("S6").Value = 3.55
$ws.Range("W6").Value = 2.04
$ws.Range("X6").Value = 13

# Row 7
$ws.Range("G7").Value = 1.85
$ws.Range("H7").Value = 4.4
$ws.Range("N7").Value = 1.01
$ws.Range("O7").Value = 1.3
$ws.Range("Q7").Value = 1.78
$ws.Range("R7").Value = 1.18
$ws.Range("S7").Value = 1.01
$ws.Range("W7").Value = 2.16
$ws.Range("Y7").Value = 980
$ws.Range("Z7").Value = 60
$ws.Range("AC7").Value = 12
$ws.Range("AD7").Value = 980
$ws.Range("AH7").Value = 980
$ws.Range("AJ7").Value = 980
$ws.Range("AL7").Value = 980

# Row 8
$ws.Range("F8").Value = 2.34
$ws.Range("G8").Value = 2.54
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 3.5
$ws.Range("J8").Value = 3.3
$ws.Range("K8").Value = 3.55
$ws.Range("L8").Value = 1.01
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 1.01
$ws.Range("O8").Value = 1.39
$ws.Range("P8").Value = 1.74
$ws.Range("R8").Value = 1.23
$ws.Range("S8").Value = 3.55
$ws.Range("T8").Value = 1.01
$ws.Range("U8").Value = 1.01
$ws.Range("V8").Value = 1.4
$ws.Range("W8").Value = 1.64
$ws.Range("X8").Value = 17
$ws.Range("Y8").Value = 17
$ws.Range("Z8").Value = 34
$ws.Range("AA8").Value = 95
$ws.Range("AB8").Value = 13
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 21
$ws.Range("AE8").Value = 65
$ws.Range("AF8").Value = 22
$ws.Range("AG8").Value = 17.5
$ws.Range("AH8").Value = 28
$ws.Range("AI8").Value = 85
$ws.Range("AJ8").Value = 50
$ws.Range("AK8").Value = 44
$ws.Range("AL8").Value = 70
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 1000
$ws.Range("AO8").Value = 1000

# Row 9
$ws.Range("G9").Value = 2.16
$ws.Range("J9").Value = 3.35
$ws.Range("L9").Value = 1.01
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 1.01
$ws.Range("O9").Value = 1.38
$ws.Range("P9").Value = 1.71
$ws.Range("Q9").Value = 1.38
$ws.Range("R9").Value = 1.24
$ws.Range("S9").Value = 3.4
$ws.Range("T9").Value = 1.01
$ws.Range("U9").Value = 1.01
$ws.Range("V9").Value = 1.21
$ws.Range("W9").Value = 1.86
$ws.Range("X9").Value = 1000
$ws.Range("Y9").Value = 1000
$ws.Range("Z9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 1000
$ws.Range("AC9").Value = 1000
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AH9").Value = 1000
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("AO9").Value = 1000
